# Refresh the cryptos price/volume snapshot (GitHub Actions data update).
# Updates Price (column D) and Volume(1h) (column E) text cells for rows 2-51
# on the single data sheet, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.968.31"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.635.60"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'212.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'0.524"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'23.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.867.42"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "1.631.99"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'65.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "27.970.93"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "'231.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "'10.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("D25").Value = "'154.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "1.408.80"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'1.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.57%  "
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "'0.564"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "'0.874"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'67.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "1.776.83"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "'87.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  +8.71%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -0.46%  "
